$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths -------------------------------------------------
# Column B gets a little wider
$ws.Columns.Item(2).ColumnWidth = 16.5
# Column E joins the C:D group (same width as C and D)
$ws.Columns.Item(5).ColumnWidth = 17.6
# New columns F and G get a narrower width
$ws.Columns.Item(6).ColumnWidth = 9.2
$ws.Columns.Item(7).ColumnWidth = 9.2

# --- Row 24: "Hill Climber" section title --------------------------
$ws.Range("A24").Value = "Hill Climber"

# --- Row 25: header row (same headers as row 2) ---------------------
$ws.Range("A25").Value = "Number of correlating pictures per class"
$ws.Range("B25").Value = "Number of parents"
$ws.Range("C25").Value = "Correctly Classified"
$ws.Range("D25").Value = "Incorrectly Classified"
$ws.Range("E25").Value = "Kappa"
$ws.Range("F25").Value = "RMS Error"

# --- Row 26-28: group "5" data --------------------------------------
$ws.Range("A26").Value = 5
$ws.Range("B26").Value = 1
$ws.Range("C26").NumberFormat = "0.0000%"
$ws.Range("C26").Value = 0.376883
$ws.Range("D26").NumberFormat = "0.0000%"
$ws.Range("D26").Formula = "=1-C26"
$ws.Range("E26").Value = 0.2849
$ws.Range("F26").Value = 0.3306

$ws.Range("A27").Value = 5
$ws.Range("B27").Value = 2
$ws.Range("C27").NumberFormat = "0.0000%"
$ws.Range("C27").Value = 0.761713
$ws.Range("D27").NumberFormat = "0.0000%"
$ws.Range("D27").Formula = "=1-C27"
$ws.Range("E27").Value = 0.7102
$ws.Range("F27").Value = 0.188

$ws.Range("A28").Value = 5
$ws.Range("B28").Value = 3
$ws.Range("C28").NumberFormat = "0.0000%"
$ws.Range("C28").Value = 0.760475
$ws.Range("D28").NumberFormat = "0.0000%"
$ws.Range("D28").Formula = "=1-C28"
$ws.Range("E28").Value = 0.7082
$ws.Range("F28").Value = 0.1852

# --- Row 29: blank spacer row (keeps percentage formatting) ---------
$ws.Range("C29").NumberFormat = "0.0000%"
$ws.Range("D29").NumberFormat = "0.0000%"

# --- Row 30-32: group "10" data --------------------------------------
$ws.Range("A30").Value = 10
$ws.Range("B30").Value = 1
$ws.Range("C30").NumberFormat = "0.0000%"
$ws.Range("C30").Value = 0.40516
$ws.Range("D30").NumberFormat = "0.0000%"
$ws.Range("D30").Formula = "=1-C30"
$ws.Range("E30").Value = 0.3183
$ws.Range("F30").Value = 0.3322

$ws.Range("A31").Value = 10
$ws.Range("B31").Value = 2
$ws.Range("C31").NumberFormat = "0.0000%"
$ws.Range("C31").Value = 0.79969
$ws.Range("D31").NumberFormat = "0.0000%"
$ws.Range("D31").Formula = "=1-C31"
$ws.Range("E31").Value = 0.7566
$ws.Range("F31").Value = 0.176

$ws.Range("A32").Value = 10
$ws.Range("B32").Value = 3
$ws.Range("C32").NumberFormat = "0.0000%"
$ws.Range("C32").Value = 0.795253
$ws.Range("D32").NumberFormat = "0.0000%"
$ws.Range("D32").Formula = "=1-C32"
$ws.Range("E32").Value = 0.7508
$ws.Range("F32").Value = 0.1736

# --- Row 33: blank spacer row (keeps percentage formatting) ---------
$ws.Range("C33").NumberFormat = "0.0000%"
$ws.Range("D33").NumberFormat = "0.0000%"

# --- Row 34-36: group "20" data --------------------------------------
$ws.Range("A34").Value = 20
$ws.Range("B34").Value = 1
$ws.Range("C34").NumberFormat = "0.0000%"
$ws.Range("C34").Value = 0.409288
$ws.Range("D34").NumberFormat = "0.0000%"
$ws.Range("D34").Formula = "=1-C34"
$ws.Range("E34").Value = 0.3231
$ws.Range("F34").Value = 0.3358

$ws.Range("A35").Value = 20
$ws.Range("B35").Value = 2
$ws.Range("C35").NumberFormat = "0.0000%"
$ws.Range("C35").Value = 0.817023
$ws.Range("D35").NumberFormat = "0.0000%"
$ws.Range("D35").Formula = "=1-C35"
$ws.Range("E35").Value = 0.7778
$ws.Range("F35").Value = 0.1719

$ws.Range("A36").Value = 20
$ws.Range("B36").Value = 3
$ws.Range("C36").NumberFormat = "0.0000%"
$ws.Range("C36").Value = 0.801238
$ws.Range("D36").NumberFormat = "0.0000%"
$ws.Range("D36").Formula = "=1-C36"
$ws.Range("E36").Value = 0.758
$ws.Range("F36").Value = 0.173

# --- Selection matches the new active cell ---------------------------
$ws.Range("G36").Select()
